$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "35.094.29"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.07%  "
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.893.17"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.74%  "
# Row 4
$ws.Range("E4").Value = "  -0.38%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "250.94"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.73%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.695"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.42%  "
# Row 7
$ws.Range("E7").Value = "  -0.29%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "41.22"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.04%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.351"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.66%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "52.04"
$ws.Range("D10").Style = "Normal"
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0745"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.52%  "
# Row 12
$ws.Range("E12").Value = "  -1.25%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.167.31"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.70%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "12.94"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.88%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.723"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.35%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.93"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.54%  "
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.879.70"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.47%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "35.072.00"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.05%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "73.58"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.58%  "
# Row 20
$ws.Range("E20").Value = "  +0.11%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "250.40"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.40%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.86"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.21%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.98"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.47%  "
# Row 24
$ws.Range("E24").Value = "  -0.34%  "
# Row 25
$ws.Range("E25").Value = "  +3.91%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.23"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.03%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "167.44"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.06%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.49"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.47%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.30"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.24%  "
# Row 30
$ws.Range("E30").Value = "  -3.32%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.128.66"
$ws.Range("D31").Style = "Normal"
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.28"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.15%  "
# Row 33
$ws.Range("E33").Value = "  +2.24%  "
# Row 34
$ws.Range("B34").Value = "WEMIXToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.93"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.98%  "
# Row 35
$ws.Range("B35").Value = "TrustWalletToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.58"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.09%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.20"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.07%  "
# Row 37
$ws.Range("E37").Value = "  -0.41%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.842"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -8.55%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.00"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.07%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.41"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.98%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "98.55"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.50%  "
# Row 42
$ws.Range("E42").Value = "  +1.33%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0659"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.37%  "
# Row 44
$ws.Range("E44").Value = "  -2.99%  "
# Row 45
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.38"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.28%  "
# Row 46
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.295.91"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.36%  "
# Row 47
$ws.Range("E47").Value = "  +0.18%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.75"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.42%  "
# Row 49
$ws.Range("E49").Value = "  +8.46%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.50"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.35%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.99"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.94%  "
